$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "time" column (old column F). The "percentage" column in
# old column G shifts left into F; all other columns are unaffected.
$ws.Columns.Item(6).Delete()

# Update the rows whose metrics/architecture changed with this run
$ws.Range("B8").Value = 'nnnn/nnnnn/nnnnn/nnn'
$ws.Range("C8").Value = 0.9857142857142858
$ws.Range("D8").Value = 0.9333333373069763
$ws.Range("B9").Value = 'nnnn/nnnnnn/nnnnnn/nnn'
$ws.Range("C9").Value = 0.9809523812362126
$ws.Range("D9").Value = 0.9666666686534882
$ws.Range("B11").Value = 'nnnn/nnnn/n/nnnn/nnnn/n/nnnn/nnn'
$ws.Range("C11").Value = 0.9904761910438538
$ws.Range("D11").Value = 1.0
$ws.Range("E11").Value = 1.0
$ws.Range("B12").Value = 'nnnn/nnnn/nnnn/nnn'
$ws.Range("C12").Value = 0.9904761904761905
$ws.Range("D12").Value = 0.9666666388511658
$ws.Range("B13").Value = 'nnnn/nnnn/nnnn/nnnn/nnn'
$ws.Range("C13").Value = 0.9809523809523809
$ws.Range("D13").Value = 0.9333333373069763
$ws.Range("E13").Value = 0.9333333373069763
$ws.Range("B14").Value = 'nnnn/nnnn/nnnn/nnnnnn/nnn'
$ws.Range("C14").Value = 0.9809523809523809
$ws.Range("D14").Value = 0.9666666388511658
$ws.Range("B15").Value = 'nnnn/nnnnn/n/nnnnn/nnnnn/n/nnnnn/nnn'
$ws.Range("D15").Value = 0.9333333373069763
$ws.Range("E15").Value = 0.9333333373069763
$ws.Range("B16").Value = 'nnnn/nnnnnn/nnnnnn/nnnnn/nnn'
$ws.Range("C16").Value = 0.9809523809523809
$ws.Range("E16").Value = 1.0
$ws.Range("B17").Value = 'nnnn/nnnnnnn/nnnnnnn/nnnnnnn/nnnnnnn/nnn'
$ws.Range("C17").Value = 0.9619047624724252
$ws.Range("D17").Value = 1.0
$ws.Range("B18").Value = 'nnnn/nnnnnnnn/nnnnnnnn/nnn'
$ws.Range("C18").Value = 0.9904761904761905
$ws.Range("D18").Value = 0.9666666388511658
$ws.Range("B19").Value = 'nnnn/nnnnnnnn/nnnnnnnn/nnnnnnnn/nnnnnnnn/nnn'
$ws.Range("D19").Value = 1.0
$ws.Range("B21").Value = 'nnnn/nnnnnnnnnnnn/nnn'
$ws.Range("C21").Value = 0.9809523815200443
